$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column F (date/time values) needs an explicit width of 15 characters.
$ws.Columns("F").ColumnWidth = 14 + 1/6

# Apply date/time number formats to column F, which holds date-time serial
# values. The first block (rows 1-4) and the repeated block (rows 7-10)
# get the same treatment:
#   * row 1 / row 7  -> F cell already carries a font/fill/border/protection
#                        style; it now additionally gets a short-date format.
#   * row 2 / row 8  -> F cell had no style; it now gets a short-date format.
#   * row 3-4 / 9-10 -> F cells had no style; they now get a time format.

$ws.Range("F1").NumberFormat = "mm-dd-yy"
$ws.Range("F2").NumberFormat = "mm-dd-yy"
$ws.Range("F3").NumberFormat = "h:mm:ss"
$ws.Range("F4").NumberFormat = "h:mm:ss"

$ws.Range("F7").NumberFormat = "mm-dd-yy"
$ws.Range("F8").NumberFormat = "mm-dd-yy"
$ws.Range("F9").NumberFormat = "h:mm:ss"
$ws.Range("F10").NumberFormat = "h:mm:ss"
